# DAS-677 - CCRU - Creation of Scenes for SOVI SOCVI
#
# The "Panoramic Photo" and "Panoramic photo of Cooler" scene tags are being
# extended with their corresponding "SS_" (subscription-service) scene
# names, so every cell that used to just say "Panoramic Photo" now reads
# "Panoramic Photo, SS_Panoramic Photo" (and similarly for the cooler photo
# scene).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HoReCa Restaurant_Cafe")
$ws.Activate()

# Column AA ("Scenes to include") on the Availability KPI rows: replace the
# old single scene tag with the combined Panoramic Photo / SS_Panoramic
# Photo tag.
$panoramicRows = @(4,5,6,7,8,9,10,12,13,15,16,17,18,19,20,21,22,24)
foreach ($r in $panoramicRows) {
    $ws.Cells.Item($r, 27).Value = "Panoramic Photo, SS_Panoramic Photo"
}

# Row 36 ("Cooler fullness" KPI), column Z ("Scenes to exclude"): replace the
# old single scene tag with the combined Panoramic photo of Cooler /
# SS_Panoramic photo of Cooler - Horeca tag.
$ws.Cells.Item(36, 26).Value = "Panoramic photo of Cooler, SS_Panoramic photo of Cooler - Horeca"

# Re-point the active selection to where the user was last working
# (selection on Z37 in the frozen bottom-left pane), leaving the existing
# freeze of row 1 untouched.
$ws.Range("Z37").Select()
